# Update "DSM Scheduled Flights vs actual.xlsx":
#   - append 40 new daily rows (2022-03-24 .. 2022-05-02) below the
#     existing data table on "Ark1"
#   - extend the D-column "% on time" shared formula down through the
#     new rows
#   - move the active selection/scroll position to reflect where the
#     user ended up after typing the new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 718
$lastNewRow  = 757
$lastOldRow  = 717

$dates = @(
    "2022-03-24","2022-03-25","2022-03-26","2022-03-27","2022-03-28",
    "2022-03-29","2022-03-30","2022-03-31","2022-04-01","2022-04-02",
    "2022-04-03","2022-04-04","2022-04-05","2022-04-06","2022-04-07",
    "2022-04-08","2022-04-09","2022-04-10","2022-04-11","2022-04-12",
    "2022-04-13","2022-04-14","2022-04-15","2022-04-16","2022-04-17",
    "2022-04-18","2022-04-19","2022-04-20","2022-04-21","2022-04-22",
    "2022-04-23","2022-04-24","2022-04-25","2022-04-26","2022-04-27",
    "2022-04-28","2022-04-29","2022-04-30","2022-05-01","2022-05-02"
)

$scheduled = @(
    71,72,60,73,69,58,61,76,77,51,62,67,61,73,87,69,49,65,63,53,
    69,78,68,52,52,67,74,72,85,68,49,60,72,63,67,70,70,59,63,69
)

$tracked = @(
    70,71,58,72,68,57,61,73,75,49,59,66,59,70,78,65,45,64,63,53,
    65,76,66,52,52,65,72,71,83,66,47,59,66,62,63,68,69,56,57,68
)

# Carry the existing row's number formats / styles (text date in A,
# whole numbers in B/C, percentage in D) down onto the new rows before
# filling them in, same as dragging the fill handle would.
$ws.Range("A$lastOldRow`:D$lastOldRow").Copy() | Out-Null
$ws.Range("A$firstNewRow`:D$lastNewRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $firstNewRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $scheduled[$i]
    $ws.Cells.Item($r, 3).Value = $tracked[$i]
}

# Fill the ratio formula down the new rows in one shot (same as
# dragging the D717 fill handle to D757) so it is written back out as
# a single shared formula, matching the rest of the column.
$ws.Range("D$firstNewRow`:D$lastNewRow").Formula = "=C$firstNewRow/B$firstNewRow"

# Reflect the final scroll/selection state from the edit.
$ws.Range("F754").Select()
